$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 51 and 52: the two match records swap places (everything except the
# shared "opening" timestamp columns K, O, S, which stay put because both
# matches opened at the same snapshot time).
# ---------------------------------------------------------------------------
$row51 = @("Brabrand","1","Skive","0",2.69,"30/09/2023 12:03",2.87,3.13,"30/09/2023 13:53",3.18,2.42,"30/09/2023 13:53",2.45,"https://www.betexplorer.com/football/denmark/2nd-division/brabrand-skive/IZaweirL/")
$row52 = @("Roskilde","3","Esbjerg","3",3.56,"30/09/2023 13:57",3.81,3.63,"30/09/2023 12:01",3.76,1.79,"30/09/2023 13:57",1.84,"https://www.betexplorer.com/football/denmark/2nd-division/roskilde-esbjerg/zBAYeXSE/")

$ws.Range("F52").Value = $row51[0]
$ws.Range("G52").Value = $row51[1]
$ws.Range("H52").Value = $row51[2]
$ws.Range("I52").Value = $row51[3]
$ws.Range("J52").Value = $row51[4]
$ws.Range("M52").Value = $row51[5]
$ws.Range("L52").Value = $row51[6]
$ws.Range("N52").Value = $row51[7]
$ws.Range("Q52").Value = $row51[8]
$ws.Range("P52").Value = $row51[9]
$ws.Range("R52").Value = $row51[10]
$ws.Range("U52").Value = $row51[11]
$ws.Range("T52").Value = $row51[12]
$ws.Range("V52").Value = $row51[13]

$ws.Range("F51").Value = $row52[0]
$ws.Range("G51").Value = $row52[1]
$ws.Range("H51").Value = $row52[2]
$ws.Range("I51").Value = $row52[3]
$ws.Range("J51").Value = $row52[4]
$ws.Range("M51").Value = $row52[5]
$ws.Range("L51").Value = $row52[6]
$ws.Range("N51").Value = $row52[7]
$ws.Range("Q51").Value = $row52[8]
$ws.Range("P51").Value = $row52[9]
$ws.Range("R51").Value = $row52[10]
$ws.Range("U51").Value = $row52[11]
$ws.Range("T51").Value = $row52[12]
$ws.Range("V51").Value = $row52[13]

# ---------------------------------------------------------------------------
# Rows 68, 69, 70: a 3-way cyclic rotation (row68 -> row69, row69 -> row70,
# row70 -> row68), timestamps K/O/S stay fixed per row.
# ---------------------------------------------------------------------------
$ws.Range("F69").Value = "Brabrand"
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = "AB Copenhagen"
$ws.Range("I69").Value = 2
$ws.Range("J69").Value = 3.31
$ws.Range("M69").Value = "21/10/2023 13:41"
$ws.Range("L69").Value = 3.55
$ws.Range("N69").Value = 3.44
$ws.Range("Q69").Value = "21/10/2023 13:41"
$ws.Range("P69").Value = 3.49
$ws.Range("R69").Value = 1.95
$ws.Range("U69").Value = "21/10/2023 13:41"
$ws.Range("T69").Value = 1.99
$ws.Range("V69").Value = "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-ab-copenhagen/Iys5yoc6/"

$ws.Range("F70").Value = "Esbjerg"
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = "FA 2000"
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 1.21
$ws.Range("M70").Value = "21/10/2023 06:38"
$ws.Range("L70").Value = 1.18
$ws.Range("N70").Value = 6.19
$ws.Range("Q70").Value = "21/10/2023 13:25"
$ws.Range("P70").Value = 7.26
$ws.Range("R70").Value = 7.96
$ws.Range("U70").Value = "21/10/2023 13:25"
$ws.Range("T70").Value = 11.24
$ws.Range("V70").Value = "https://www.betexplorer.com/football/denmark/2nd-division/esbjerg-frederiksberg-alliancen-2000/00W9z5CC/"

$ws.Range("F68").Value = "Skive"
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = "Roskilde"
$ws.Range("I68").Value = 2
$ws.Range("J68").Value = 3.5
$ws.Range("M68").Value = "21/10/2023 13:41"
$ws.Range("L68").Value = 3.25
$ws.Range("N68").Value = 3.38
$ws.Range("Q68").Value = "21/10/2023 13:41"
$ws.Range("P68").Value = 3.43
$ws.Range("R68").Value = 1.91
$ws.Range("U68").Value = "21/10/2023 13:41"
$ws.Range("T68").Value = 2.12
$ws.Range("V68").Value = "https://www.betexplorer.com/football/denmark/2nd-division/skive-roskilde/vor1xRs0/"

# ---------------------------------------------------------------------------
# Rows 75 and 76: swap places, same pattern as rows 51/52.
# ---------------------------------------------------------------------------
$ws.Range("F75").Value = "Aarhus Fremad"
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = "Middelfart"
$ws.Range("I75").Value = 2
$ws.Range("J75").Value = 1.73
$ws.Range("M75").Value = "28/10/2023 13:51"
$ws.Range("L75").Value = 1.82
$ws.Range("N75").Value = 3.76
$ws.Range("Q75").Value = "28/10/2023 13:51"
$ws.Range("P75").Value = 3.73
$ws.Range("R75").Value = 3.69
$ws.Range("U75").Value = "28/10/2023 13:51"
$ws.Range("T75").Value = 3.94
$ws.Range("V75").Value = "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-middelfart/Iy3ziQdg/"

$ws.Range("F76").Value = "Roskilde"
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = "Brabrand"
$ws.Range("I76").Value = 4
$ws.Range("J76").Value = 1.56
$ws.Range("M76").Value = "27/10/2023 10:36"
$ws.Range("L76").Value = 1.48
$ws.Range("N76").Value = 3.92
$ws.Range("Q76").Value = "28/10/2023 12:03"
$ws.Range("P76").Value = 4.27
$ws.Range("R76").Value = 4.56
$ws.Range("U76").Value = "27/10/2023 10:36"
$ws.Range("T76").Value = 6.09
$ws.Range("V76").Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-brabrand/Wv7vjpCa/"

# ---------------------------------------------------------------------------
# New row 79: one more match appended at the end of the table.
# Copy formatting from the last existing data row first, then fill values.
# ---------------------------------------------------------------------------
$ws.Range("A78:V78").Copy()
$ws.Range("A79:V79").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "denmark"
$ws.Range("C79").Value = "2nd-division"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45233.79166666666
$ws.Range("F79").Value = "Skive"
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = "Thisted FC"
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 2.55
$ws.Range("K79").Value = "02/11/2023 07:12"
$ws.Range("L79").Value = 2.23
$ws.Range("M79").Value = "03/11/2023 18:58"
$ws.Range("N79").Value = 3.17
$ws.Range("O79").Value = "02/11/2023 07:12"
$ws.Range("P79").Value = 3.48
$ws.Range("Q79").Value = "03/11/2023 18:58"
$ws.Range("R79").Value = 2.46
$ws.Range("S79").Value = "02/11/2023 07:12"
$ws.Range("T79").Value = 2.99
$ws.Range("U79").Value = "03/11/2023 18:56"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/denmark/2nd-division/skive-thisted-fc/xpIKVFQp/"
